# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values: B=TB, C=d2S, D=K, E=IP, F=Win (unchanged), G=sum
$data = @{
    2 = @(0.7287194209349384, 1.65323645889881,  157.8057217802531,  6.48142807727062,  166.6691057373575)
    3 = @(1.505614041169197, 0.3375848360084654, 157.8057217802531,  6.48142807727062,  166.1303487347014)
    4 = @(0.7287194209349384, 0.05231270169004087, 3.082599426703578, 6.48142807727062,  10.34505962659918)
    5 = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    6 = @(3.182878228561681, 1.65323645889881,  0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    7 = @(3.182878228561681, 1.65323645889881,  16.98373111632243,  0.4998867070740569, 22.31973251085698)
    8 = @(0.7287194209349384, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 1.433824611717217)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($r, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($r, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($r, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($r, 7).Value = $vals[4]   # G - sum
}

$wb.Save()
